# Refresh crypto market data (prices & 1h volume change) pulled from coinranking.
# Also reflects a rank swap between Bittensor and VeChain (rows 45-46).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.795.66'
$ws.Range("E2").Value = '  -1.59%  '
$ws.Range("D3").Value = '2.904.14'
$ws.Range("E3").Value = '  -2.56%  '
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '589.78'
$ws.Range("E5").Value = '  +0.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.57'
$ws.Range("E6").Value = '  -2.82%  '
$ws.Range("E7").Value = '  -0.20%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.505'
$ws.Range("E8").Value = '  +0.27%  '
$ws.Range("D9").Value = '2.904.80'
$ws.Range("E9").Value = '  -2.14%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.78'
$ws.Range("E10").Value = '  +1.12%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.144'
$ws.Range("E11").Value = '  -1.53%  '
$ws.Range("E12").Value = '  -2.46%  '
$ws.Range("E13").Value = '  +0.07%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.48'
$ws.Range("E14").Value = '  -3.33%  '
$ws.Range("E15").Value = '  +0.51%  '
$ws.Range("D16").Value = '3.385.02'
$ws.Range("E16").Value = '  -2.63%  '
$ws.Range("D17").Value = '60.775.69'
$ws.Range("E17").Value = '  -1.82%  '
$ws.Range("E18").Value = '  -3.20%  '
$ws.Range("D19").Value = '2.902.45'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '429.42'
$ws.Range("E20").Value = '  -2.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.49'
$ws.Range("E21").Value = '  -3.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.683'
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("E23").Value = '  -4.32%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '80.98'
$ws.Range("E24").Value = '  -0.60%  '
$ws.Range("E25").Value = '  -2.48%  '
$ws.Range("E26").Value = '  -0.22%  '
$ws.Range("E27").Value = '  -0.10%  '
$ws.Range("E28").Value = '  -0.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.29'
$ws.Range("E29").Value = '  +3.17%  '
$ws.Range("E31").Value = '  -1.54%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.05'
$ws.Range("E32").Value = '  -3.76%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.58'
$ws.Range("E33").Value = '  -1.76%  '
$ws.Range("E34").Value = '  -0.72%  '
$ws.Range("D35").Value = '0.0₃0849'
$ws.Range("E35").Value = '  +1.14%  '
$ws.Range("E36").Value = '  -1.23%  '
$ws.Range("E37").Value = '  -2.69%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.04'
$ws.Range("E38").Value = '  +1.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '49.53'
$ws.Range("E39").Value = '  -0.98%  '
$ws.Range("E40").Value = '  -2.16%  '
$ws.Range("E41").Value = '  -2.26%  '
$ws.Range("E42").Value = '  -2.57%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.287'
$ws.Range("E43").Value = '  -3.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.09'
$ws.Range("E44").Value = '  -5.99%  '
$ws.Range("B45").Value = 'Bittensor'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '374.78'
$ws.Range("E45").Value = '  +0.95%  '
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0347'
$ws.Range("E46").Value = '  -1.00%  '
$ws.Range("D47").Value = '2.702.13'
$ws.Range("E47").Value = '  +1.07%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '129.17'
$ws.Range("E48").Value = '  -2.15%  '
$ws.Range("E49").Value = '  -0.01%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.99'
$ws.Range("E50").Value = '  -6.90%  '
$ws.Range("E51").Value = '  -0.80%  '
